# Natmi LR-pair output (Sema3g-Nrp2) rebuilt "following Dr Hou advice":
# the sending/target cluster grid now uses all 3 clusters (ECs, FAPs, sCs)
# sending to all 3 clusters (a full 3x3 = 9 row table instead of the old
# 2-cluster x 3-cluster = 6 row table), with refreshed numeric results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T for data rows 2..10 (row 1 is the header, left untouched).
$data = @(
    ,@("ECs",  "Sema3g", "Nrp2", "ECs",  2, 0.6666666666666666, 3.949943,          11.849829,         0.7136716036992452,  0.7136716036992451,  3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 105.9802004267263,  953.8218038405371, 0.4222549990529236,  0.4222549990529236)
    ,@("ECs",  "Sema3g", "Nrp2", "FAPs", 2, 0.6666666666666666, 3.949943,          11.849829,         0.7136716036992452,  0.7136716036992451,  3, 1, 10.21969166666667, 30.659075,          0.2253617819930474, 0.2253617819930474, 40.36719956090833,  363.304796048175,  0.1608343043674978,  0.1608343043674978)
    ,@("ECs",  "Sema3g", "Nrp2", "sCs",  2, 0.6666666666666666, 3.949943,          11.849829,         0.7136716036992452,  0.7136716036992451,  3, 1, 8.297426666666667, 24.89228,           0.1829725319067811, 0.1829725319067811, 32.77436238001333,  294.96926142012,   0.1305823002788238,  0.1305823002788237)
    ,@("FAPs", "Sema3g", "Nrp2", "ECs",  3, 1,                  0.625675,          1.877025,          0.1130463099453651,  0.1130463099453651,  3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 16.78737184359167,  151.086346592325,  0.06688562253491709, 0.06688562253491709)
    ,@("FAPs", "Sema3g", "Nrp2", "FAPs", 3, 1,                  0.625675,          1.877025,          0.1130463099453651,  0.1130463099453651,  3, 1, 10.21969166666667, 30.659075,          0.2253617819930474, 0.2253617819930474, 6.394205583541667,  57.547850251875,   0.02547631785702583, 0.02547631785702583)
    ,@("FAPs", "Sema3g", "Nrp2", "sCs",  3, 1,                  0.625675,          1.877025,          0.1130463099453651,  0.1130463099453651,  3, 1, 8.297426666666667, 24.89228,           0.1829725319067811, 0.1829725319067811, 5.191492429666666,  46.723431867,      0.02068436955342217, 0.02068436955342217)
    ,@("sCs",  "Sema3g", "Nrp2", "ECs",  3, 1,                  0.9590606666666667,2.877182,          0.1732820863553897,  0.1732820863553897,  3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 25.73238187860512,  231.591436907446,  0.1025250645123309,  0.1025250645123308)
    ,@("sCs",  "Sema3g", "Nrp2", "FAPs", 3, 1,                  0.9590606666666667,2.877182,          0.1732820863553897,  0.1732820863553897,  3, 1, 10.21969166666667, 30.659075,          0.2253617819930474, 0.2253617819930474, 9.801304302961112,  88.21173872665,    0.03905115976852376, 0.03905115976852375)
    ,@("sCs",  "Sema3g", "Nrp2", "sCs",  3, 1,                  0.9590606666666667,2.877182,          0.1732820863553897,  0.1732820863553897,  3, 1, 8.297426666666667, 24.89228,           0.1829725319067811, 0.1829725319067811, 7.957735550551111,  71.61961995496,    0.03170586207453514, 0.03170586207453514)
)

# Write each value individually (cell-by-cell) so existing rows 2-7 are
# overwritten in place and new rows 8-10 are appended, growing the sheet's
# used range from A1:T7 to A1:T10.
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
